$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("C8").Value = 7
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "1792.00"

# Row 9
$ws.Range("C9").Value = 96
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "4"
$ws.Range("E9").Value = "Long point  (up to 10 mtr.)"
$ws.Range("F9").Value = 662
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "63552.00"

# Row 10
$ws.Range("A10").Value = ""
$ws.Range("C10").Value = 75
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.0"
$ws.Range("E10").Value = "Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it's  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet's & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F10").Value = 0
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "0.00"

# Row 11
$ws.Range("C11").Value = 31
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.0"
$ws.Range("E11").Value = "P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F11").Value = 50
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "1550.00"

# Row 12
$ws.Range("A12").Value = "R. mtr."
$ws.Range("C12").Value = 42
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17"
$ws.Range("E12").Value = "25 mm"
$ws.Range("F12").Value = 56
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "2352.00"

# Row 13
$ws.Range("A13").Value = ""
$ws.Range("C13").Value = 4
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.0"
$ws.Range("E13").Value = "Supplying and drawing FR PVC insulated & unsheathed flexible copper conductor as per PWD specification for electrical Works with ISI marked (IS:694) and as per IS 8130 : 2013 of 1.1 kV grade . Wire should be made from  99.90 % purity copper, class 2 stranding in acc. to IS:8130/IEC 60228 for  lower watt loss , oxygen free for less chances of oxidization, insulation PVC type A/C/D , flame retardant as per IS 10810-53, better amperage rating as per IS:3961 part 5,  in existing  surface or recessed PVC/ MS conduit/casing capping making connections with Copper Lugs of suitable size, Ferrules,testing etc. as required. OEM Must have its own in house NABL lab setup for all testing facilities for wires.   For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F13").Value = 0
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "0.00"

# Row 14
$ws.Range("A14").Value = "Mtr."
$ws.Range("C14").Value = 71
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19"
$ws.Range("E14").Value = "2 x 2.5 sq. mm. + 1x1.5sqmm"
$ws.Range("F14").Value = 81
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "5751.00"

# Row 15
$ws.Range("A15").Value = "Set"
$ws.Range("C15").Value = 48
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.0"
$ws.Range("E15").Value = "Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. 'B' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR .   "
$ws.Range("F15").Value = 5733
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "275184.00"

# Row 16
$ws.Range("A16").Value = ""
$ws.Range("C16").Value = 27
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.0"
$ws.Range("E16").Value = "Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F16").Value = 0
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "0.00"

# Row 17
$ws.Range("A17").Value = "Each"
$ws.Range("C17").Value = 71
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25"
$ws.Range("E17").Value = "1200 mm Sweep BEE 1 Star rated (service value >=4.0 to < 4.5 )"
$ws.Range("F17").Value = 1890
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "134190.00"

# Row 18
$ws.Range("A18").Value = ""
$ws.Range("C18").Value = 12
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.0"
$ws.Range("E18").Value = "Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F18").Value = 0
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "0.00"

# Row 19
$ws.Range("A19").Value = "Each"
$ws.Range("C19").Value = 87
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "27"
$ws.Range("E19").Value = "1170mm(+/-10%) LED batten with min. lumen output 2200 lm"
$ws.Range("F19").Value = 492
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "42804.00"

# Row 20
$ws.Range("C20").Value = 45
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.0"
$ws.Range("E20").Value = "Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"

# Row 21
$ws.Range("A21").Value = ""
$ws.Range("C21").Value = 51
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29"
$ws.Range("E21").Value = "Single pole MCB   (With B/C curve tripping Characteristics)"
$ws.Range("F21").Value = 0
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "0.00"

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "31"
$ws.Range("E22").Value = "Double pole MCB(With B/C curve tripping Characteristics)"

# Row 23
$ws.Range("A23").Value = "Each"
$ws.Range("C23").Value = 89
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "32"
$ws.Range("E23").Value = " 50/63 A rating"
$ws.Range("F23").Value = 900
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "80100.00"

# Row 24
$ws.Range("C24").Value = 38
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "18.0"
$ws.Range("E24").Value = "Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"

# Row 25
$ws.Range("A25").Value = ""
$ws.Range("C25").Value = 12
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "36"
$ws.Range("E25").Value = "Total"
$ws.Range("F25").Value = 0
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "0.00"

# Row 26
$ws.Range("A26").Value = "%"
$ws.Range("C26").Value = 42
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37"
$ws.Range("E26").Value = "Add Tender Premium "

# Row 27
$ws.Range("A27").Value = ""
$ws.Range("B27:I27").ClearContents()

# Row 28
$ws.Range("B28").Value = ""
$ws.Range("C28").Value = ""
$ws.Range("D28").Value = ""
$ws.Range("E28").Value = "Grand Total Rs."
$ws.Range("F28").Value = ""
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "607275.00"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "607275.00"

# Row 29
$ws.Range("A29").Value = ""
$ws.Range("B29").Value = ""
$ws.Range("C29").Value = ""
$ws.Range("D29").Value = ""
$ws.Range("E29").Value = "Tender Premium @ 0%"
$ws.Range("F29").Value = ""
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "0.00"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "0.00"

# Row 30
$ws.Range("B30").Value = ""
$ws.Range("C30").Value = ""
$ws.Range("D30").Value = ""
$ws.Range("E30").Value = "NET PAYABLE AMOUNT Rs."
$ws.Range("F30").Value = ""
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "607275.00"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "607275.00"

# Remove now-unused trailing rows 31-34 (content moved up / consolidated)
$ws.Rows("31:34").Delete()
